# Auto-generated Excel COM-interop script
# Applies the Leve profit recalculations recorded in the commit diff,
# sheet by sheet, row by row, matching the before/after cell values.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 124: Luncheon Bound / Luncheon Toadskin Codex
$ws.Range("H124").Value = 34866.668
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 34866.668
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 34866.668
$ws.Range("M124").ClearContents()
$ws.Range("N124").Value = -44686.668

# Row 126: Rebuilding to Code / Saigaskin Codex
$ws.Range("H126").Value = 41523.332
$ws.Range("J126").Value = 41523.332
$ws.Range("L126").Value = 41523.332
$ws.Range("N126").Value = -51403.332

# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 2885702
$ws.Range("I137").Value = 1163828.4
$ws.Range("K137").Value = 3491485.2
$ws.Range("M137").Value = -3488935.2

# Row 140: Tome for Tradition / Book of Ra'Kaznar
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# Row 43: They've Got Legs / Steel Sabatons
$ws.Range("H43").Value = 9663.4
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 9663.4
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 9663.4
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -10289.4

# Row 44: Very Slow Array / Mythril Plate
$ws.Range("H44").Value = 34949
$ws.Range("J44").Value = 34949
$ws.Range("L44").Value = 34949
$ws.Range("N44").Value = -35925

# Row 55: Employee Retention / Mythril Elmo
$ws.Range("H55").Value = 16016.333
$ws.Range("J55").Value = 16016.333
$ws.Range("L55").Value = 16016.333
$ws.Range("N55").Value = -16646.333

# Row 80: A Squire to Inspire / Titanium Hoplon
$ws.Range("H80").Value = 22560.572
$ws.Range("J80").Value = 22560.572
$ws.Range("L80").Value = 22560.572
$ws.Range("N80").Value = -24556.572

# Row 83: All's Fair in Highborn Assassination (L) / Titanium Hoplon
$ws.Range("H83").Value = 22560.572
$ws.Range("J83").Value = 22560.572
$ws.Range("L83").Value = 67681.716
$ws.Range("N83").Value = -77665.716

# Row 127: Once and for Alchemy / Bismuth Alembic
$ws.Range("H127").Value = 46986.668
$ws.Range("J127").Value = 46986.668
$ws.Range("L127").Value = 46986.668
$ws.Range("N127").Value = -56906.668

$ws = $wb.Worksheets.Item("BSM")
# Row 35: Lancers' Creed / Crowsbeak Hammer
$ws.Range("H35").Value = 34864.4
$ws.Range("J35").Value = 34864.4
$ws.Range("L35").Value = 34864.4
$ws.Range("N35").Value = -35484.4

# Row 82: Spirituality Inspector / Titanium Lump Hammer
$ws.Range("H82").Value = 52146.176
$ws.Range("I82").Value = 109171.336
$ws.Range("J82").Value = 32019.646
$ws.Range("K82").Value = 109171.336
$ws.Range("L82").Value = 32019.646
$ws.Range("M82").Value = -108788.336
$ws.Range("N82").Value = -32785.646

# Row 85: The Clamor for Hammers (L) / Titanium Lump Hammer
$ws.Range("H85").Value = 52146.176
$ws.Range("I85").Value = 109171.336
$ws.Range("J85").Value = 32019.646
$ws.Range("K85").Value = 109171.336
$ws.Range("L85").Value = 32019.646
$ws.Range("M85").Value = -107845.336
$ws.Range("N85").Value = -34671.646

# Row 99: Meddle in Metal / Oroshigane Ingot
$ws.Range("H99").Value = 1428.0952
$ws.Range("I99").Value = 805
$ws.Range("K99").Value = 805
$ws.Range("M99").Value = 693

# Row 122: To Delight a Dancer / High Durium Tathlums
$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

# Row 124: History of the Hrothgar / High Durium Bayonet
$ws.Range("H124").Value = 40389.5
$ws.Range("J124").Value = 40389.5
$ws.Range("L124").Value = 40389.5
$ws.Range("N124").Value = -50209.5

# Row 126: Records of the Republic / Bismuth War Scythe
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()

# Row 135: Axes to the Maxes / Ruthenium War Axe
$ws.Range("H135").Value = 33231
$ws.Range("J135").Value = 33231
$ws.Range("L135").Value = 33231
$ws.Range("N135").Value = -43371

$ws = $wb.Worksheets.Item("CRP")
# Row 20: Re-crating the Scene / Iron Spear
$ws.Range("H20").Value = 48299.715
$ws.Range("J20").Value = 48299.715
$ws.Range("L20").Value = 48299.715
$ws.Range("N20").Value = -48771.715

# Row 30: Polearms Aplenty / Iron Spear
$ws.Range("H30").Value = 48299.715
$ws.Range("J30").Value = 48299.715
$ws.Range("L30").Value = 48299.715
$ws.Range("N30").Value = -48481.715

# Row 41: The Lone Bowman / Oak Longbow
$ws.Range("H41").Value = 17899
$ws.Range("I41").Value = 4900
$ws.Range("J41").Value = 21148.75
$ws.Range("K41").Value = 4900
$ws.Range("L41").Value = 21148.75
$ws.Range("M41").Value = -4472
$ws.Range("N41").Value = -22004.75

# Row 50: The Arsenal of Theocracy / Cobalt Halberd
$ws.Range("H50").Value = 19588.334
$ws.Range("J50").Value = 19588.334
$ws.Range("L50").Value = 19588.334
$ws.Range("N50").Value = -20838.334

# Row 51: Greenstone for Greenhorns / Jade Crook
$ws.Range("H51").Value = 17932.1
$ws.Range("J51").Value = 19591.223
$ws.Range("L51").Value = 19591.223
$ws.Range("N51").Value = -21063.223

# Row 60: Bowing to Greater Power / Yew Longbow
$ws.Range("H60").Value = 22767
$ws.Range("J60").Value = 22767
$ws.Range("L60").Value = 22767
$ws.Range("N60").Value = -23789

# Row 61: Incant Now, Think Later / Jade Crook
$ws.Range("H61").Value = 17932.1
$ws.Range("J61").Value = 19591.223
$ws.Range("L61").Value = 19591.223
$ws.Range("N61").Value = -20287.223

# Row 62: Splinter in the Sewers / Cedar Lumber
$ws.Range("H62").Value = 4000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 4000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 4000
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -5248

# Row 65: The Lumber of Their Discontent (L) / Cedar Lumber
$ws.Range("H65").Value = 4000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 4000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 20000
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -26240

# Row 68: Do You Even String Bow / Holy Cedar Composite Bow
$ws.Range("H68").Value = 19295
$ws.Range("J68").Value = 19295
$ws.Range("L68").Value = 19295
$ws.Range("N68").Value = -20793

# Row 71: Win One Bow, Get Three Free (L) / Holy Cedar Composite Bow
$ws.Range("H71").Value = 19295
$ws.Range("J71").Value = 19295
$ws.Range("L71").Value = 57885
$ws.Range("N71").Value = -65373

# Row 109: Playing the Market / White Oak Necklace
$ws.Range("H109").Value = 10833.333
$ws.Range("J109").Value = 11000
$ws.Range("L109").Value = 11000
$ws.Range("N109").Value = -13080

# Row 127: In Rod We Trust / Red Pine Fishing Rod
$ws.Range("H127").Value = 54380
$ws.Range("J127").Value = 54380
$ws.Range("L127").Value = 54380
$ws.Range("N127").Value = -64300

# Row 128: An A-prop-riate Request / Ironwood Spear
$ws.Range("H128").Value = 48299.715
$ws.Range("J128").Value = 48299.715
$ws.Range("L128").Value = 48299.715
$ws.Range("N128").Value = -58259.715

# Row 130: Annals of the Empire II / Integral Magitek Rod
$ws.Range("H130").Value = 59980
$ws.Range("J130").Value = 59980
$ws.Range("L130").Value = 59980
$ws.Range("N130").Value = -70020

$ws = $wb.Worksheets.Item("GSM")
# Row 93: One Ring Circus / Triphane Ring of Slaying
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()

# Row 122: Awarding Academic Excellence / Ametrine
$ws.Range("H122").Value = 3342.611
$ws.Range("I122").Value = 4123.375
$ws.Range("J122").Value = 2718
$ws.Range("K122").Value = 12370.125
$ws.Range("L122").Value = 8154
$ws.Range("M122").Value = -9920.125
$ws.Range("N122").Value = -13054

# Row 123: Workplace Workout / Ametrine Ring of Fending
$ws.Range("H123").Value = 34108.4
$ws.Range("J123").Value = 34108.4
$ws.Range("L123").Value = 34108.4
$ws.Range("N123").Value = -39008.4

$ws = $wb.Worksheets.Item("LTW")
# Row 92: Vested Interest / Gaganaskin Vest
$ws.Range("H92").Value = 38788.5
$ws.Range("J92").Value = 38788.5
$ws.Range("L92").Value = 38788.5
$ws.Range("N92").Value = -43780.5

# Row 109: Band Substances / Smilodonskin Wristband
$ws.Range("H109").Value = 24583.5
$ws.Range("J109").Value = 24583.5
$ws.Range("L109").Value = 24583.5
$ws.Range("N109").Value = -27357.5

# Row 122: Hell on Leather / Gaja Leather
$ws.Range("H122").Value = 3026.7917
$ws.Range("I122").Value = 2339.4546
$ws.Range("J122").Value = 3608.3845
$ws.Range("K122").Value = 7018.3638
$ws.Range("L122").Value = 10825.1535
$ws.Range("M122").Value = -4568.3638
$ws.Range("N122").Value = -15725.1535

# Row 127: Loyal Turncoat / Saigaskin Coat of Fending
$ws.Range("H127").Value = 46330
$ws.Range("J127").Value = 46330
$ws.Range("L127").Value = 46330
$ws.Range("N127").Value = -56250

$ws = $wb.Worksheets.Item("WVR")
# Row 93: What Guides Want / Bloodhempen Doublet of Crafting
$ws.Range("H93").Value = 27561.166
$ws.Range("J93").Value = 27561.166
$ws.Range("L93").Value = 27561.166
$ws.Range("N93").Value = -32553.166

# Row 109: Turban in Training / Brightlinen Turban of Crafting
$ws.Range("H109").Value = 9977.777
$ws.Range("J109").Value = 9977.777
$ws.Range("L109").Value = 9977.777
$ws.Range("N109").Value = -12751.777

# Row 113: A Tender Table / Pixie Floss
$ws.Range("H113").Value = 39453.69
$ws.Range("I113").Value = 59146.824
$ws.Range("J113").Value = 2255.5557
$ws.Range("K113").Value = 177440.472
$ws.Range("L113").Value = 6766.6671
$ws.Range("M113").Value = -175270.472
$ws.Range("N113").Value = -11106.6671

# Row 125: Color Coated / Almasty Serge Coat of Healing
$ws.Range("H125").Value = 35782.918
$ws.Range("J125").Value = 35782.918
$ws.Range("L125").Value = 35782.918
$ws.Range("N125").Value = -45622.918

# Row 127: Turban Sprawl / Snow Linen Turban of Crafting
$ws.Range("H127").Value = 35344.375
$ws.Range("J127").Value = 35344.375
$ws.Range("L127").Value = 35344.375
$ws.Range("N127").Value = -45264.375
